# "using gains for all"
# Insert two new data columns (M_TotalTax, M_CorpTax) right after the
# M_POP column (E) and before the old "GFA - Sales" column (old F),
# shifting the existing GFA/IMF/OECD columns two places to the right
# (old F:M -> new H:O), then populate the two new columns with the
# "gains" figures for every group row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank columns at F:G - this shifts old F:M (GFA - Sales ...
# OECD - Sales + Emp) to H:O automatically, matching the dimension
# growing from A1:M11 to A1:O11.
$ws.Range("F1:G1").EntireColumn.Insert()

# New header labels
$ws.Range("F1").Value = "M_TotalTax"
$ws.Range("G1").Value = "M_CorpTax"

# New data values per group row (EU, G20, G20_noOECD, G24, G7, G77,
# G7_noUS, OECD, OECD_noUS, US)
$ws.Range("F2").Value = 6308727034979.312
$ws.Range("G2").Value = 399825921028.5854

$ws.Range("F3").Value = 16630145391623.02
$ws.Range("G3").Value = 1639742485782.957

$ws.Range("F4").Value = 4450994137606.095
$ws.Range("G4").Value = 601350231413.5104

$ws.Range("F5").Value = 4183547438952.192
$ws.Range("G5").Value = 598849276038.3025

$ws.Range("F6").Value = 11223287075501.79
$ws.Range("G6").Value = 872292028558.4308

$ws.Range("F7").Value = 1841737275230.086
$ws.Range("G7").Value = 214321200777.9413

$ws.Range("F8").Value = 6192585801479.285
$ws.Range("G8").Value = 516695167857.3162

$ws.Range("F9").Value = 14653861967257.56
$ws.Range("G9").Value = 1232540278767.842

$ws.Range("F10").Value = 9623160693235.053
$ws.Range("G10").Value = 876943418066.7275

$ws.Range("F11").Value = 5030701274022.499
$ws.Range("G11").Value = 355596860701.1148
